$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - values updated
$ws.Range("B3").Value = 0.9880307178679533
$ws.Range("C3").Value = 0.9885354706119558
$ws.Range("D3").Value = 0.8011762619519383

# Row 4: model name and values updated
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.986212995291484
$ws.Range("C4").Value = 0.9868506764944788
$ws.Range("D4").Value = 0.7869406469453467

# Row 5: model name and values updated
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.860458804105096
$ws.Range("C5").Value = 0.8466266397326364
$ws.Range("D5").Value = 0.5964040945023409
